$p = $ppt.ActivePresentation

$notesText = @{
    1 = "早上好，同学们！在我们之前的讲座中，我们学习了牛顿运动定律。今天，我们将通过引入功、能量和功率来扩展这些概念。这些概念对于理解力如何在距离和时间上影响物体至关重要。"
    2 = "在物理学中，功被定义为力引起位移。如果一个力使物体发生位移，则该力对物体做了功。做功的大小计算为力与位移的乘积，具体来说，是力沿位移方向的分量与位移大小的乘积。"
    3 = "能量是做功的能力。能量有多种形式，例如动能（运动的能量）和势能（储存的能量）。能量守恒定律指出，能量既不能被创造，也不能被毁灭，只能从一种形式转化为另一种形式。"
    4 = "动能是物体由于其运动而具有的能量。它与物体的质量及其速度的平方成正比。因此，速度更快或质量更大的物体具有更大的动能。"
    5 = "势能是储存的能量。重力势能取决于物体的质量、高度和重力加速度。弹性势能储存在被拉伸或压缩的弹性材料中。"
    6 = "最后，功率是做功或能量转移的速率。它告诉我们做功的速度有多快。功率等于功除以时间，或者力乘以速度。这就结束了我们关于功、能量和功率的讲座。请为下一节课的简短测验做好准备。"
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $np = $s.NotesPage
    for ($j = 1; $j -le $np.Shapes.Count; $j++) {
        $shp = $np.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.Name -like "Notes Placeholder*") {
            $shp.TextFrame.TextRange.Text = $notesText[$i]
        }
    }
}
